$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 417.3125
$ws.Cells.Item(28, 9).Value = 467.69232
$ws.Cells.Item(28, 10).Value = 199
$ws.Cells.Item(28, 11).Value = 467.69232
$ws.Cells.Item(28, 12).Value = 199
$ws.Cells.Item(28, 13).Value = 17.30768
$ws.Cells.Item(28, 14).Value = -1169

$ws.Cells.Item(33, 8).Value = 873.7
$ws.Cells.Item(33, 9).Value = 705
$ws.Cells.Item(33, 10).Value = 1267.3334
$ws.Cells.Item(33, 11).Value = 705
$ws.Cells.Item(33, 12).Value = 1267.3334
$ws.Cells.Item(33, 13).Value = -476
$ws.Cells.Item(33, 14).Value = -1725.3334

$ws.Cells.Item(69, 8).Value = 3835.75
$ws.Cells.Item(69, 10).Value = 3910
$ws.Cells.Item(69, 12).Value = 11730
$ws.Cells.Item(69, 14).Value = -13478

$ws.Cells.Item(72, 8).Value = 3835.75
$ws.Cells.Item(72, 10).Value = 3910
$ws.Cells.Item(72, 12).Value = 35190
$ws.Cells.Item(72, 14).Value = -43926

$ws.Cells.Item(101, 8).Value = 1309.6923
$ws.Cells.Item(101, 9).Value = 293.27274
$ws.Cells.Item(101, 10).Value = 6900
$ws.Cells.Item(101, 11).Value = 879.81822
$ws.Cells.Item(101, 12).Value = 20700
$ws.Cells.Item(101, 13).Value = 742.18178
$ws.Cells.Item(101, 14).Value = -23944

$ws.Cells.Item(113, 8).Value = 4960.3335
$ws.Cells.Item(113, 9).Value = 3394.75
$ws.Cells.Item(113, 10).Value = 5743.125
$ws.Cells.Item(113, 11).Value = 3394.75
$ws.Cells.Item(113, 12).Value = 5743.125
$ws.Cells.Item(113, 13).Value = -140.75
$ws.Cells.Item(113, 14).Value = -12251.125

$ws.Cells.Item(132, 8).Value = 6253971.5
$ws.Cells.Item(132, 9).Value = 7410603.5
$ws.Cells.Item(132, 10).Value = 8158.8
$ws.Cells.Item(132, 11).Value = 22231810.5
$ws.Cells.Item(132, 12).Value = 24476.4
$ws.Cells.Item(132, 13).Value = -22229280.5
$ws.Cells.Item(132, 14).Value = -29536.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2876.0952
$ws.Cells.Item(61, 9).Value = 2024.4286
$ws.Cells.Item(61, 10).Value = 4579.4287
$ws.Cells.Item(61, 11).Value = 2024.4286
$ws.Cells.Item(61, 12).Value = 4579.4287
$ws.Cells.Item(61, 13).Value = -1812.4286
$ws.Cells.Item(61, 14).Value = -5003.4287

$ws.Cells.Item(102, 8).Value = 4643.4165
$ws.Cells.Item(102, 9).Value = 2724.875
$ws.Cells.Item(102, 11).Value = 2724.875
$ws.Cells.Item(102, 13).Value = -1102.875

$ws.Cells.Item(110, 8).Value = 1254.7
$ws.Cells.Item(110, 9).Value = 529.5
$ws.Cells.Item(110, 10).Value = 2946.8333
$ws.Cells.Item(110, 11).Value = 529.5
$ws.Cells.Item(110, 12).Value = 2946.8333
$ws.Cells.Item(110, 13).Value = 1515.5
$ws.Cells.Item(110, 14).Value = -7036.8333

$ws.Cells.Item(122, 8).Value = 2542.5454
$ws.Cells.Item(122, 9).Value = 1971.0714
$ws.Cells.Item(122, 10).Value = 5742.8
$ws.Cells.Item(122, 11).Value = 5913.2142
$ws.Cells.Item(122, 12).Value = 17228.4
$ws.Cells.Item(122, 13).Value = -3463.2142
$ws.Cells.Item(122, 14).Value = -22128.4

$ws.Cells.Item(136, 8).Value = 2876.0952
$ws.Cells.Item(136, 9).Value = 2024.4286
$ws.Cells.Item(136, 10).Value = 4579.4287
$ws.Cells.Item(136, 11).Value = 6073.2858
$ws.Cells.Item(136, 12).Value = 13738.2861
$ws.Cells.Item(136, 13).Value = -3523.2858
$ws.Cells.Item(136, 14).Value = -18838.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1837.2903
$ws.Cells.Item(20, 9).Value = 1389.3636
$ws.Cells.Item(20, 11).Value = 1389.3636
$ws.Cells.Item(20, 13).Value = -1142.3636

$ws.Cells.Item(94, 8).Value = 1034.6923
$ws.Cells.Item(94, 9).Value = 516.25
$ws.Cells.Item(94, 10).Value = 1265.1111
$ws.Cells.Item(94, 11).Value = 516.25
$ws.Cells.Item(94, 12).Value = 1265.1111
$ws.Cells.Item(94, 13).Value = -65.25
$ws.Cells.Item(94, 14).Value = -2167.1111

$ws.Cells.Item(99, 8).Value = 3463.2
$ws.Cells.Item(99, 9).Value = 2944.9
$ws.Cells.Item(99, 10).Value = 4499.8
$ws.Cells.Item(99, 11).Value = 2944.9
$ws.Cells.Item(99, 12).Value = 4499.8
$ws.Cells.Item(99, 13).Value = -1446.9
$ws.Cells.Item(99, 14).Value = -7495.8

$ws.Cells.Item(113, 8).Value = 5000
$ws.Cells.Item(113, 9).Value = 5000
$ws.Cells.Item(113, 11).Value = 5000
$ws.Cells.Item(113, 13).Value = -2830

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 843.4737
$ws.Cells.Item(22, 9).Value = 277.75
$ws.Cells.Item(22, 11).Value = 277.75
$ws.Cells.Item(22, 13).Value = 72.25

$ws.Cells.Item(31, 8).Value = 2175.2424
$ws.Cells.Item(31, 9).Value = 1287.4286
$ws.Cells.Item(31, 10).Value = 2829.4211
$ws.Cells.Item(31, 11).Value = 1287.4286
$ws.Cells.Item(31, 12).Value = 2829.4211
$ws.Cells.Item(31, 13).Value = -992.4286
$ws.Cells.Item(31, 14).Value = -3419.4211

$ws.Cells.Item(34, 8).Value = 2175.2424
$ws.Cells.Item(34, 9).Value = 1287.4286
$ws.Cells.Item(34, 10).Value = 2829.4211
$ws.Cells.Item(34, 11).Value = 1287.4286
$ws.Cells.Item(34, 12).Value = 2829.4211
$ws.Cells.Item(34, 13).Value = -1085.4286
$ws.Cells.Item(34, 14).Value = -3233.4211

$ws.Cells.Item(107, 8).Value = 1432.3077
$ws.Cells.Item(107, 9).Value = 820.3333
$ws.Cells.Item(107, 10).Value = 1956.8572
$ws.Cells.Item(107, 11).Value = 820.3333
$ws.Cells.Item(107, 12).Value = 1956.8572
$ws.Cells.Item(107, 13).Value = 1099.6667
$ws.Cells.Item(107, 14).Value = -5796.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 2863.6365
$ws.Cells.Item(80, 9).Value = 2166.6667
$ws.Cells.Item(80, 10).Value = 3125
$ws.Cells.Item(80, 11).Value = 6500.000100000001
$ws.Cells.Item(80, 12).Value = 9375
$ws.Cells.Item(80, 13).Value = -5564.000100000001
$ws.Cells.Item(80, 14).Value = -11247

$ws.Cells.Item(83, 8).Value = 2863.6365
$ws.Cells.Item(83, 9).Value = 2166.6667
$ws.Cells.Item(83, 10).Value = 3125
$ws.Cells.Item(83, 11).Value = 19500.0003
$ws.Cells.Item(83, 12).Value = 28125
$ws.Cells.Item(83, 13).Value = -14820.0003
$ws.Cells.Item(83, 14).Value = -37485

$ws.Cells.Item(87, 8).Value = 11462.5
$ws.Cells.Item(87, 9).Value = 6975
$ws.Cells.Item(87, 11).Value = 20925
$ws.Cells.Item(87, 13).Value = -19677

$ws.Cells.Item(90, 8).Value = 11462.5
$ws.Cells.Item(90, 9).Value = 6975
$ws.Cells.Item(90, 11).Value = 62775
$ws.Cells.Item(90, 13).Value = -56535

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 12982.625
$ws.Cells.Item(57, 10).Value = 12965.25
$ws.Cells.Item(57, 12).Value = 12965.25
$ws.Cells.Item(57, 14).Value = -14605.25

$ws.Cells.Item(70, 8).Value = 4467.857
$ws.Cells.Item(70, 9).Value = 4800
$ws.Cells.Item(70, 10).Value = 4218.75
$ws.Cells.Item(70, 11).Value = 4800
$ws.Cells.Item(70, 12).Value = 4218.75
$ws.Cells.Item(70, 13).Value = -4530
$ws.Cells.Item(70, 14).Value = -4758.75

$ws.Cells.Item(73, 8).Value = 4467.857
$ws.Cells.Item(73, 9).Value = 4800
$ws.Cells.Item(73, 10).Value = 4218.75
$ws.Cells.Item(73, 11).Value = 4800
$ws.Cells.Item(73, 12).Value = 4218.75
$ws.Cells.Item(73, 13).Value = -3864
$ws.Cells.Item(73, 14).Value = -6090.75

$ws.Cells.Item(102, 8).Value = 46831.78
$ws.Cells.Item(102, 9).Value = 2452.625
$ws.Cells.Item(102, 11).Value = 2452.625
$ws.Cells.Item(102, 13).Value = -830.625

$ws.Cells.Item(113, 8).Value = 2292.7144
$ws.Cells.Item(113, 9).Value = 812.25
$ws.Cells.Item(113, 10).Value = 4266.6665
$ws.Cells.Item(113, 11).Value = 812.25
$ws.Cells.Item(113, 12).Value = 4266.6665
$ws.Cells.Item(113, 13).Value = 1357.75
$ws.Cells.Item(113, 14).Value = -8606.666499999999

$ws.Cells.Item(122, 8).Value = 4508.227
$ws.Cells.Item(122, 9).Value = 3821.6924
$ws.Cells.Item(122, 10).Value = 5499.8887
$ws.Cells.Item(122, 11).Value = 11465.0772
$ws.Cells.Item(122, 12).Value = 16499.6661
$ws.Cells.Item(122, 13).Value = -9015.0772
$ws.Cells.Item(122, 14).Value = -21399.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 166673840
$ws.Cells.Item(61, 9).Value = 500004500
$ws.Cells.Item(61, 10).Value = 8498.75
$ws.Cells.Item(61, 11).Value = 500004500
$ws.Cells.Item(61, 12).Value = 8498.75
$ws.Cells.Item(61, 13).Value = -500004298
$ws.Cells.Item(61, 14).Value = -8902.75

$ws.Cells.Item(82, 8).Value = 3745.2
$ws.Cells.Item(82, 9).Value = 3207.4285
$ws.Cells.Item(82, 10).Value = 5000
$ws.Cells.Item(82, 11).Value = 3207.4285
$ws.Cells.Item(82, 12).Value = 5000
$ws.Cells.Item(82, 13).Value = -2846.4285
$ws.Cells.Item(82, 14).Value = -5722

$ws.Cells.Item(85, 8).Value = 3745.2
$ws.Cells.Item(85, 9).Value = 3207.4285
$ws.Cells.Item(85, 10).Value = 5000
$ws.Cells.Item(85, 11).Value = 3207.4285
$ws.Cells.Item(85, 12).Value = 5000
$ws.Cells.Item(85, 13).Value = -1959.4285
$ws.Cells.Item(85, 14).Value = -7496

$ws.Cells.Item(113, 8).Value = 166673840
$ws.Cells.Item(113, 9).Value = 500004500
$ws.Cells.Item(113, 10).Value = 8498.75
$ws.Cells.Item(113, 11).Value = 500004500
$ws.Cells.Item(113, 12).Value = 8498.75
$ws.Cells.Item(113, 13).Value = -500002330
$ws.Cells.Item(113, 14).Value = -12838.75

$ws.Cells.Item(131, 8).Value = 30004.334
$ws.Cells.Item(131, 10).Value = 30004.334
$ws.Cells.Item(131, 12).Value = 30004.334
$ws.Cells.Item(131, 14).Value = -40084.334

$ws.Cells.Item(132, 8).Value = 2929.9714
$ws.Cells.Item(132, 9).Value = 2055.2104
$ws.Cells.Item(132, 10).Value = 3968.75
$ws.Cells.Item(132, 11).Value = 6165.6312
$ws.Cells.Item(132, 12).Value = 11906.25
$ws.Cells.Item(132, 13).Value = -3635.6312
$ws.Cells.Item(132, 14).Value = -16966.25

$ws.Cells.Item(136, 8).Value = 3034357.8
$ws.Cells.Item(136, 9).Value = 4548650
$ws.Cells.Item(136, 10).Value = 5773.1816
$ws.Cells.Item(136, 11).Value = 13645950
$ws.Cells.Item(136, 12).Value = 17319.5448
$ws.Cells.Item(136, 13).Value = -13643400
$ws.Cells.Item(136, 14).Value = -22419.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 244556.1
$ws.Cells.Item(132, 9).Value = 347373.4
$ws.Cells.Item(132, 10).Value = 31577.357
$ws.Cells.Item(132, 11).Value = 1042120.2
$ws.Cells.Item(132, 12).Value = 94732.071
$ws.Cells.Item(132, 13).Value = -1039590.2
$ws.Cells.Item(132, 14).Value = -99792.071

$ws.Cells.Item(136, 8).Value = 1982.9318
$ws.Cells.Item(136, 9).Value = 1549.9032
$ws.Cells.Item(136, 10).Value = 3015.5386
$ws.Cells.Item(136, 11).Value = 4649.7096
$ws.Cells.Item(136, 12).Value = 9046.6158
$ws.Cells.Item(136, 13).Value = -2099.7096
$ws.Cells.Item(136, 14).Value = -14146.6158
